$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-23T08:28:04+00:00"

# --- 2. "Mapping Table 1" sheet: remap rows 3-15 (columns A and D) ---
# Source column A switches from the FRLMImageIllustrative business-model
# names to the FRCDAImageIllustrative CDA names; target column D switches
# from FRCDAImageIllustrative CDA names to FRMediaDocument FHIR names.
$ws = $wb.Worksheets.Item("Mapping Table 1")

$rows = @(
    @{ Row = 3;  A = "FRCDAImageIllustrative.id";                  D = "FRMediaDocument.identifier" },
    @{ Row = 4;  A = "FRCDAImageIllustrative.languageCode";        D = "FRMediaDocument.content.language" },
    @{ Row = 5;  A = "FRCDAImageIllustrative.value";               D = "FRMediaDocument.content.data" },
    @{ Row = 6;  A = "FRCDAImageIllustrative.value.mediaType";     D = "FRMediaDocument.content.contentType" },
    @{ Row = 7;  A = "FRCDAImageIllustrative.subject";             D = "FRMediaDocument.subject:Patient" },
    @{ Row = 8;  A = "FRCDAImageIllustrative.specimen";            D = "FRMediaDocument.subject:Specimen" },
    @{ Row = 9;  A = "FRCDAImageIllustrative.performer";           D = "FRMediaDocument.operator.extension:performer" },
    @{ Row = 10; A = "FRCDAImageIllustrative.author";              D = "FRMediaDocument.operator.extension:author" },
    @{ Row = 11; A = "FRCDAImageIllustrative.informant";           D = "FRMediaDocument.operator.extension:informant" },
    @{ Row = 12; A = "FRCDAImageIllustrative.participant";         D = "FRMediaDocument.operator.extension:participant" },
    @{ Row = 13; A = "FRCDAImageIllustrative.entryRelationship";   D = "FRMediaDocument.basedOn" },
    @{ Row = 14; A = "FRCDAImageIllustrative.reference";           D = "FRMediaDocument.partOf" },
    @{ Row = 15; A = "FRCDAImageIllustrative.precondition";        D = "FRMediaDocument.reasonCode" }
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 4).Value = $item.D
}
